# DangVH: Commit update tien do du an
# Update the Coding-status (I) and UI-merge-status (J) for the 5 "Admin
# account" use cases (UC041-UC045, rows 49-53 on Sheet1) that DangVH owns
# (column F = "DangVH"): they move from "Not Start" to "In Progress" (I)
# and "Not Start" to "Done" (J). The TONG HOP summary sheet recalculates
# its COUNTIF tallies automatically from these cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Match the existing colour-coded "cell style" convention used throughout
# the sheet: In Progress -> Neutral (yellow), Done -> Good (green). Copy
# the format from cells that already carry the right look (keeps borders
# intact, unlike assigning .Style directly).
$ws.Range("I48").Copy() | Out-Null
$ws.Range("I49:I53").PasteSpecial(-4122) | Out-Null

$ws.Range("I7").Copy() | Out-Null
$ws.Range("J49:J53").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

$ws.Range("I49:I53").Value = "In Progress"
$ws.Range("J49:J53").Value = "Done"

# Leave the selection where the editor ended up.
$ws.Range("E73").Select() | Out-Null

Write-Output "Updated status for UC041-UC045 (DangVH)"
